$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "22.438.67"
$ws.Range("E2").Value = "  +0.21%  "

# Row 3
$ws.Range("D3").Value = "1.574.18"
$ws.Range("E3").Value = "  +0.32%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("E5").Value = "  +0.20%  "

# Row 6
$ws.Range("D6").Value = "'291.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "

# Row 7
$ws.Range("D7").Value = "'0.3768"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.62%  "

# Row 8
$ws.Range("D8").Value = "'49.91"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.64%  "

# Row 9
$ws.Range("E9").Value = "  +1.52%  "

# Row 10
$ws.Range("D10").Value = "'1.163"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.72%  "

# Row 11
$ws.Range("D11").Value = "'0.07672"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.18%  "

# Row 12
$ws.Range("E12").Value = "  +0.20%  "

# Row 13
$ws.Range("D13").Value = "'21.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.87%  "

# Row 14
$ws.Range("D14").Value = "'6.011"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.66%  "

# Row 15
$ws.Range("E15").Value = "  +1.21%  "

# Row 16
$ws.Range("D16").Value = "1.573.36"
$ws.Range("E16").Value = "  +0.22%  "

# Row 17
$ws.Range("D17").Value = "'0.00001135"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.86%  "

# Row 18
$ws.Range("D18").Value = "'90.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.34%  "

# Row 19
$ws.Range("D19").Value = "'0.06758"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.67%  "

# Row 20
$ws.Range("E20").Value = "  +0.16%  "

# Row 21
$ws.Range("D21").Value = "'16.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.54%  "

# Row 22
$ws.Range("E22").Value = "  -0.40%  "

# Row 23
$ws.Range("E23").Value = "  +0.70%  "

# Row 24
$ws.Range("B24").Value = "WrappedBTC"
$ws.Range("C24").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D24").Value = "22.444.02"
$ws.Range("E24").Value = "  +0.19%  "

# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.428"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.92%  "

# Row 26
$ws.Range("D26").Value = "'2.737"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.03%  "

# Row 27
$ws.Range("D27").Value = "'20.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.32%  "

# Row 28
$ws.Range("D28").Value = "'146.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.07%  "

# Row 29
$ws.Range("D29").Value = "'5.034"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.86%  "

# Row 30
$ws.Range("D30").Value = "'126.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.08%  "

# Row 31
$ws.Range("D31").Value = "1.749.87"
$ws.Range("E31").Value = "  +0.23%  "

# Row 32
$ws.Range("D32").Value = "'6.209"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.80%  "

# Row 33
$ws.Range("D33").Value = "'2.016"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.94%  "

# Row 34
$ws.Range("E34").Value = "  +1.33%  "

# Row 35
$ws.Range("D35").Value = "'10.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.05%  "

# Row 36
$ws.Range("D36").Value = "'0.08577"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.88%  "

# Row 37
$ws.Range("D37").Value = "'0.02556"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.99%  "

# Row 38
$ws.Range("E38").Value = "  +0.89%  "

# Row 39
$ws.Range("D39").Value = "'0.06583"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.28%  "

# Row 40
$ws.Range("E40").Value = "  +6.88%  "

# Row 41
$ws.Range("E41").Value = "  -0.78%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.6461"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.20%  "

# Row 43
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'11.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.89%  "

# Row 44
$ws.Range("D44").Value = "'14.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.81%  "

# Row 45
$ws.Range("E45").Value = "  +0.16%  "

# Row 46
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "'3.801"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.79%  "

# Row 47
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.6022"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.09%  "

# Row 48
$ws.Range("E48").Value = "  +9.18%  "

# Row 49
$ws.Range("E49").Value = "  -1.57%  "

# Row 50
$ws.Range("D50").Value = "'125.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.19%  "

# Row 51
$ws.Range("D51").Value = "'0.07330"
$ws.Range("D51").Style = "Normal"
